# feat!: removal of option `fieldMatchType`
# New default is `labelTypeBrackets`: header labels now carry their mapped
# field name in square brackets (e.g. "ID" -> "ID[product_ID]") so that
# there is no ambiguity between a column's display label and its bound type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 holds the column headers; append the "[fieldName]" suffix to each.
$ws.Range("A1").Value = "ID[product_ID]"
$ws.Range("B1").Value = "Quantity[quantity]"
$ws.Range("C1").Value = "ProductTitle[title]"
$ws.Range("D1").Value = "UnitPrice[price]"
$ws.Range("E1").Value = "validFrom[validFrom]"
$ws.Range("F1").Value = "timestamp[timestamp]"
$ws.Range("G1").Value = "date[date]"
$ws.Range("H1").Value = "time[time]"
$ws.Range("I1").Value = "WRONGCOLUMN[TEST]"

# Update the sheet's saved selection to span the header row (A1:H1).
$ws.Range("A1:H1").Select() | Out-Null
